$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so that
# numeric-looking values (e.g. "302.91") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.185.98'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '2.321.90'
$ws.Range("E3").Value = '  +0.85%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '302.91'
$ws.Range("D6").Value = '99.49'
$ws.Range("E6").Value = '  +0.30%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").Value = '35.96'
$ws.Range("E10").Value = '  +4.61%  '
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").Value = '17.65'
$ws.Range("E13").Value = '  -2.71%  '
$ws.Range("D14").Value = '6.92'
$ws.Range("E14").Value = '  +1.74%  '
$ws.Range("D15").Value = '2.683.93'
$ws.Range("E15").Value = '  +0.92%  '
$ws.Range("D16").Value = '2.317.22'
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("D17").Value = '0.797'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").Value = '43.095.46'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("E19").Value = '  +6.63%  '
$ws.Range("D20").Value = '6.25'
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("E21").Value = '  +0.78%  '
$ws.Range("D22").Value = '68.12'
$ws.Range("E22").Value = '  +0.32%  '
$ws.Range("D23").Value = '240.28'
$ws.Range("E23").Value = '  +1.60%  '
$ws.Range("D24").Value = '2.16'
$ws.Range("E24").Value = '  -2.16%  '
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("D27").Value = '25.51'
$ws.Range("E27").Value = '  +3.00%  '
$ws.Range("D28").Value = '168.02'
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("D29").Value = '34.25'
$ws.Range("E29").Value = '  +1.36%  '
$ws.Range("D30").Value = '9.20'
$ws.Range("E30").Value = '  +0.46%  '
$ws.Range("E31").Value = '  -1.98%  '
$ws.Range("E32").Value = '  +2.89%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = '4.74'
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("D35").Value = '17.69'
$ws.Range("E35").Value = '  +4.39%  '
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("D37").Value = '0.0698'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("E38").Value = '  +2.01%  '
$ws.Range("E39").Value = '  +0.15%  '
$ws.Range("E40").Value = '  -1.19%  '
$ws.Range("D42").Value = '1.992.99'
$ws.Range("E42").Value = '  -0.23%  '
$ws.Range("E43").Value = '  +1.32%  '
$ws.Range("D44").Value = '2.24'
$ws.Range("E44").Value = '  -4.66%  '
$ws.Range("D45").Value = '10.10'
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("D46").Value = '17.57'
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("D48").Value = '76.33'
$ws.Range("E48").Value = '  +8.82%  '
$ws.Range("D49").Value = '55.04'
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("D50").Value = '2.86'
$ws.Range("E50").Value = '  +12.86%  '
$ws.Range("D51").Value = '2.548.58'
